# Apply crypto price/volume updates from the Nov 2 2023 GitHub Actions refresh.
# Every text value is written with a leading apostrophe so Excel keeps it as
# literal text (matching the original inline-string cells) instead of silently
# re-interpreting numeric-looking strings (e.g. '1.60' or '0.1000') as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.846.10"
$ws.Range("E2").Value = "'  +0.87%  "

$ws.Range("D3").Value = "'1.811.59"
$ws.Range("E3").Value = "'  +0.55%  "

$ws.Range("E4").Value = "'  +0.45%  "

$ws.Range("D5").Value = "'231.75"
$ws.Range("E5").Value = "'  +3.27%  "

$ws.Range("D6").Value = "'0.603"
$ws.Range("E6").Value = "'  +0.47%  "

$ws.Range("E7").Value = "'  +0.41%  "

$ws.Range("D8").Value = "'40.09"
$ws.Range("E8").Value = "'  -3.58%  "

$ws.Range("D9").Value = "'0.306"
$ws.Range("E9").Value = "'  +4.97%  "

$ws.Range("D10").Value = "'0.0681"
$ws.Range("E10").Value = "'  +2.30%  "

$ws.Range("D11").Value = "'0.1000"
$ws.Range("E11").Value = "'  +0.30%  "

$ws.Range("D12").Value = "'2.072.41"

$ws.Range("D13").Value = "'1.810.88"
$ws.Range("E13").Value = "'  +0.69%  "

$ws.Range("E14").Value = "'  +1.09%  "

$ws.Range("D15").Value = "'0.656"
$ws.Range("E15").Value = "'  +4.40%  "

$ws.Range("E16").Value = "'  +5.51%  "

$ws.Range("D17").Value = "'34.823.15"
$ws.Range("E17").Value = "'  +1.03%  "

$ws.Range("D18").Value = "'68.79"
$ws.Range("E18").Value = "'  +2.30%  "

$ws.Range("E19").Value = "'  +1.99%  "

$ws.Range("D20").Value = "'236.51"
$ws.Range("E20").Value = "'  -1.54%  "

$ws.Range("D21").Value = "'11.73"
$ws.Range("E21").Value = "'  +5.34%  "

$ws.Range("E22").Value = "'  +9.62%  "

$ws.Range("E23").Value = "'  +0.45%  "

$ws.Range("E24").Value = "'  +4.14%  "

$ws.Range("D25").Value = "'172.81"
$ws.Range("E25").Value = "'  +0.54%  "

$ws.Range("D26").Value = "'7.74"
$ws.Range("E26").Value = "'  +1.04%  "

$ws.Range("D27").Value = "'17.32"
$ws.Range("E27").Value = "'  -0.30%  "

$ws.Range("E28").Value = "'  -0.89%  "

$ws.Range("D29").Value = "'1.60"
$ws.Range("E29").Value = "'  +30.44%  "

$ws.Range("E30").Value = "'  +0.57%  "

$ws.Range("D31").Value = "'3.340.14"
$ws.Range("E31").Value = "'  +37.47%  "

$ws.Range("D32").Value = "'0.0544"
$ws.Range("E32").Value = "'  +6.16%  "

$ws.Range("D33").Value = "'3.86"
$ws.Range("E33").Value = "'  +1.76%  "

$ws.Range("E34").Value = "'  +1.98%  "

$ws.Range("D35").Value = "'1.77"
$ws.Range("E35").Value = "'  -1.08%  "

$ws.Range("D36").Value = "'93.56"
$ws.Range("E36").Value = "'  +6.87%  "

$ws.Range("D37").Value = "'1.12"
$ws.Range("E37").Value = "'  +6.94%  "

$ws.Range("D38").Value = "'0.674"
$ws.Range("E38").Value = "'  +4.20%  "

$ws.Range("D39").Value = "'1.306.00"
$ws.Range("E39").Value = "'  -0.91%  "

$ws.Range("E40").Value = "'  +4.27%  "

$ws.Range("E41").Value = "'  +1.65%  "

$ws.Range("D42").Value = "'14.86"
$ws.Range("E42").Value = "'  +0.78%  "

$ws.Range("D43").Value = "'0.981"
$ws.Range("E43").Value = "'  +4.94%  "

$ws.Range("D44").Value = "'2.32"
$ws.Range("E44").Value = "'  -0.66%  "

$ws.Range("E45").Value = "'  +0.39%  "

$ws.Range("D46").Value = "'2.75"
$ws.Range("E46").Value = "'  -1.43%  "

$ws.Range("D47").Value = "'6.17"
$ws.Range("E47").Value = "'  +6.80%  "

$ws.Range("D48").Value = "'0.0514"
$ws.Range("E48").Value = "'  -0.82%  "

$ws.Range("D49").Value = "'1.986.94"
$ws.Range("E49").Value = "'  +1.20%  "

$ws.Range("E50").Value = "'  +0.44%  "

$ws.Range("B51").Value = "'Cronos"
$ws.Range("C51").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.0638"
$ws.Range("E51").Value = "'  +5.00%  "
